$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 235 (was row 236's data)
$ws.Cells.Item(235, 2).Value = 6870268
$ws.Cells.Item(235, 5).Value = 'Petrolul Ploiesti'
$ws.Cells.Item(235, 6).Value = 'ACS Sepsi'
$ws.Cells.Item(235, 7).Value = 1
$ws.Cells.Item(235, 8).Value = 2
$ws.Cells.Item(235, 10).Value = 1
$ws.Cells.Item(235, 11).Value = 'A'
$ws.Cells.Item(235, 12).Value = 2.8
$ws.Cells.Item(235, 13).Value = 3
$ws.Cells.Item(235, 14).Value = 2.55
$ws.Cells.Item(235, 15).Value = 3
$ws.Cells.Item(235, 16).Value = 3.2
$ws.Cells.Item(235, 17).Value = 2.3
$ws.Cells.Item(235, 19).Value = 1.85
$ws.Cells.Item(235, 20).Value = 2
$ws.Cells.Item(235, 22).Value = 1.875
$ws.Cells.Item(235, 23).Value = 1.975
$ws.Cells.Item(235, 25).Value = -1
$ws.Cells.Item(235, 26).Value = 1.3
$ws.Cells.Item(235, 27).Value = -1
$ws.Cells.Item(235, 28).Value = 1
$ws.Cells.Item(235, 29).Value = 0.875
$ws.Cells.Item(235, 30).Value = -1

# Row 236 (was row 235's data)
$ws.Cells.Item(236, 2).Value = 6865915
$ws.Cells.Item(236, 5).Value = 'FC Voluntari'
$ws.Cells.Item(236, 6).Value = 'Universitatea Cluj'
$ws.Cells.Item(236, 7).Value = 0
$ws.Cells.Item(236, 8).Value = 0
$ws.Cells.Item(236, 10).Value = 0
$ws.Cells.Item(236, 11).Value = 'D'
$ws.Cells.Item(236, 12).Value = 3.5
$ws.Cells.Item(236, 13).Value = 3.25
$ws.Cells.Item(236, 14).Value = 2.05
$ws.Cells.Item(236, 15).Value = 3.4
$ws.Cells.Item(236, 16).Value = 3.1
$ws.Cells.Item(236, 17).Value = 2.15
$ws.Cells.Item(236, 19).Value = 1.975
$ws.Cells.Item(236, 20).Value = 1.875
$ws.Cells.Item(236, 22).Value = 2.05
$ws.Cells.Item(236, 23).Value = 1.75
$ws.Cells.Item(236, 25).Value = 2.1
$ws.Cells.Item(236, 26).Value = -1
$ws.Cells.Item(236, 27).Value = 0.4875
$ws.Cells.Item(236, 28).Value = -0.5
$ws.Cells.Item(236, 29).Value = -1
$ws.Cells.Item(236, 30).Value = 0.75

# Row 238 (was row 239's data)
$ws.Cells.Item(238, 2).Value = 6836277
$ws.Cells.Item(238, 5).Value = 'CFR Cluj'
$ws.Cells.Item(238, 6).Value = 'AFC Hermannstadt'
$ws.Cells.Item(238, 7).Value = 1
$ws.Cells.Item(238, 11).Value = 'H'
$ws.Cells.Item(238, 12).Value = 1.7
$ws.Cells.Item(238, 14).Value = 5
$ws.Cells.Item(238, 15).Value = 1.65
$ws.Cells.Item(238, 16).Value = 3.5
$ws.Cells.Item(238, 17).Value = 5.25
$ws.Cells.Item(238, 18).Value = -0.75
$ws.Cells.Item(238, 19).Value = 1.85
$ws.Cells.Item(238, 20).Value = 2
$ws.Cells.Item(238, 21).Value = 2.25
$ws.Cells.Item(238, 22).Value = 1.875
$ws.Cells.Item(238, 23).Value = 1.975
$ws.Cells.Item(238, 24).Value = 0.6499999999999999
$ws.Cells.Item(238, 25).Value = -1
$ws.Cells.Item(238, 27).Value = 0.425
$ws.Cells.Item(238, 30).Value = 0.9750000000000001

# Row 239 (was row 238's data)
$ws.Cells.Item(239, 2).Value = 6861095
$ws.Cells.Item(239, 5).Value = 'FC Botosani'
$ws.Cells.Item(239, 6).Value = 'Farul Constanta'
$ws.Cells.Item(239, 7).Value = 0
$ws.Cells.Item(239, 11).Value = 'D'
$ws.Cells.Item(239, 12).Value = 3.75
$ws.Cells.Item(239, 14).Value = 1.909
$ws.Cells.Item(239, 15).Value = 3.1
$ws.Cells.Item(239, 16).Value = 3
$ws.Cells.Item(239, 17).Value = 2.375
$ws.Cells.Item(239, 18).Value = 0.25
$ws.Cells.Item(239, 19).Value = 1.775
$ws.Cells.Item(239, 20).Value = 2.1
$ws.Cells.Item(239, 21).Value = 2
$ws.Cells.Item(239, 22).Value = 1.8
$ws.Cells.Item(239, 23).Value = 2.05
$ws.Cells.Item(239, 24).Value = -1
$ws.Cells.Item(239, 25).Value = 2
$ws.Cells.Item(239, 27).Value = 0.3875
$ws.Cells.Item(239, 30).Value = 1.05

# Row 309 (was row 313's data)
$ws.Cells.Item(309, 2).Value = 8191462
$ws.Cells.Item(309, 5).Value = 'CSM Politehnica Iasi'
$ws.Cells.Item(309, 6).Value = 'Petrolul Ploiesti'
$ws.Cells.Item(309, 7).Value = 2
$ws.Cells.Item(309, 8).Value = 0
$ws.Cells.Item(309, 11).Value = 'H'
$ws.Cells.Item(309, 12).Value = 2.1
$ws.Cells.Item(309, 14).Value = 3.1
$ws.Cells.Item(309, 15).Value = 1.8
$ws.Cells.Item(309, 16).Value = 3.2
$ws.Cells.Item(309, 17).Value = 4.2
$ws.Cells.Item(309, 18).Value = -0.5
$ws.Cells.Item(309, 22).Value = 2.025
$ws.Cells.Item(309, 23).Value = 1.825
$ws.Cells.Item(309, 24).Value = 0.8
$ws.Cells.Item(309, 26).Value = -1
$ws.Cells.Item(309, 27).Value = 0.8500000000000001
$ws.Cells.Item(309, 28).Value = -1
$ws.Cells.Item(309, 29).Value = -0.5
$ws.Cells.Item(309, 30).Value = 0.4125

# Row 311 (was row 309's data)
$ws.Cells.Item(311, 2).Value = 8191475
$ws.Cells.Item(311, 5).Value = 'FC U Craiova 1948'
$ws.Cells.Item(311, 6).Value = 'AFC Hermannstadt'
$ws.Cells.Item(311, 7).Value = 1
$ws.Cells.Item(311, 8).Value = 3
$ws.Cells.Item(311, 9).Value = 0
$ws.Cells.Item(311, 11).Value = 'A'
$ws.Cells.Item(311, 12).Value = 2.625
$ws.Cells.Item(311, 13).Value = 3.3
$ws.Cells.Item(311, 14).Value = 2.45
$ws.Cells.Item(311, 15).Value = 2.05
$ws.Cells.Item(311, 17).Value = 3
$ws.Cells.Item(311, 18).Value = -0.25
$ws.Cells.Item(311, 22).Value = 1.825
$ws.Cells.Item(311, 23).Value = 2.025
$ws.Cells.Item(311, 24).Value = -1
$ws.Cells.Item(311, 26).Value = 2
$ws.Cells.Item(311, 27).Value = -1
$ws.Cells.Item(311, 28).Value = 1
$ws.Cells.Item(311, 29).Value = 0.825
$ws.Cells.Item(311, 30).Value = -1

# Row 312 (was row 311's data)
$ws.Cells.Item(312, 2).Value = 8191523
$ws.Cells.Item(312, 5).Value = 'Otelul Galati'
$ws.Cells.Item(312, 6).Value = 'FC Botosani'
$ws.Cells.Item(312, 12).Value = 1.666
$ws.Cells.Item(312, 13).Value = 3.6
$ws.Cells.Item(312, 14).Value = 4.6
$ws.Cells.Item(312, 15).Value = 2.9
$ws.Cells.Item(312, 16).Value = 3.5
$ws.Cells.Item(312, 17).Value = 2.2
$ws.Cells.Item(312, 18).Value = 0.25
$ws.Cells.Item(312, 19).Value = 1.85
$ws.Cells.Item(312, 20).Value = 2
$ws.Cells.Item(312, 21).Value = 2.25
$ws.Cells.Item(312, 22).Value = 1.875
$ws.Cells.Item(312, 23).Value = 1.975
$ws.Cells.Item(312, 24).Value = 1.9
$ws.Cells.Item(312, 27).Value = 0.8500000000000001
$ws.Cells.Item(312, 29).Value = -0.5
$ws.Cells.Item(312, 30).Value = 0.4875

# Row 313 (was row 312's data)
$ws.Cells.Item(313, 2).Value = 8191463
$ws.Cells.Item(313, 5).Value = 'Dinamo Bucharest'
$ws.Cells.Item(313, 6).Value = 'ACS UTA Batrana Doamna'
$ws.Cells.Item(313, 9).Value = 2
$ws.Cells.Item(313, 12).Value = 1.833
$ws.Cells.Item(313, 13).Value = 3.4
$ws.Cells.Item(313, 14).Value = 3.6
$ws.Cells.Item(313, 15).Value = 1.5
$ws.Cells.Item(313, 16).Value = 4.333
$ws.Cells.Item(313, 17).Value = 5
$ws.Cells.Item(313, 18).Value = -1
$ws.Cells.Item(313, 19).Value = 1.875
$ws.Cells.Item(313, 20).Value = 1.975
$ws.Cells.Item(313, 21).Value = 3
$ws.Cells.Item(313, 24).Value = 0.5
$ws.Cells.Item(313, 27).Value = 0.875
$ws.Cells.Item(313, 29).Value = -1
$ws.Cells.Item(313, 30).Value = 0.825

# Row 315 (was row 316's data)
$ws.Cells.Item(315, 2).Value = 7951774
$ws.Cells.Item(315, 5).Value = 'CS U Craiova'
$ws.Cells.Item(315, 6).Value = 'ACS Sepsi'
$ws.Cells.Item(315, 7).Value = 3
$ws.Cells.Item(315, 8).Value = 2
$ws.Cells.Item(315, 9).Value = 0
$ws.Cells.Item(315, 10).Value = 1
$ws.Cells.Item(315, 12).Value = 1.6
$ws.Cells.Item(315, 14).Value = 5.25
$ws.Cells.Item(315, 15).Value = 1.5
$ws.Cells.Item(315, 17).Value = 5
$ws.Cells.Item(315, 18).Value = -1
$ws.Cells.Item(315, 19).Value = 1.8
$ws.Cells.Item(315, 20).Value = 2.05
$ws.Cells.Item(315, 22).Value = 1.9
$ws.Cells.Item(315, 23).Value = 1.95
$ws.Cells.Item(315, 24).Value = 0.5
$ws.Cells.Item(315, 27).Value = 0
$ws.Cells.Item(315, 28).Value = 0
$ws.Cells.Item(315, 29).Value = 0.8999999999999999

# Row 316 (was row 315's data)
$ws.Cells.Item(316, 2).Value = 7951773
$ws.Cells.Item(316, 5).Value = 'CFR Cluj'
$ws.Cells.Item(316, 6).Value = 'Farul Constanta'
$ws.Cells.Item(316, 7).Value = 5
$ws.Cells.Item(316, 8).Value = 1
$ws.Cells.Item(316, 9).Value = 4
$ws.Cells.Item(316, 10).Value = 0
$ws.Cells.Item(316, 12).Value = 1.571
$ws.Cells.Item(316, 14).Value = 5.5
$ws.Cells.Item(316, 15).Value = 1.42
$ws.Cells.Item(316, 17).Value = 6.5
$ws.Cells.Item(316, 18).Value = -1.25
$ws.Cells.Item(316, 19).Value = 1.85
$ws.Cells.Item(316, 20).Value = 2
$ws.Cells.Item(316, 22).Value = 2.025
$ws.Cells.Item(316, 23).Value = 1.825
$ws.Cells.Item(316, 24).Value = 0.4199999999999999
$ws.Cells.Item(316, 27).Value = 0.8500000000000001
$ws.Cells.Item(316, 28).Value = -1
$ws.Cells.Item(316, 29).Value = 1.025
